$wb = $excel.ActiveWorkbook

# The "LR" execution-times sheet gets a new timing run appended as row 12
# (same Method/Computer/Language as the existing rows, a fresh Execution_time
# and a Timestamp captured later the same day).
$ws = $wb.Worksheets.Item("LR")

$ws.Range("A12").Value = "LR"
$ws.Range("B12").Value = "Windows Ryzen 9 5900x 32GB"
$ws.Range("C12").Value = 14.4005200862885
$ws.Range("D12").Value = "2025-04-12 23:46:40"
$ws.Range("E12").Value = "R"
